$wb = $excel.ActiveWorkbook

# --- Sheet "train": add span value for the Uber row, rezoom the view ---
$wsTrain = $wb.Worksheets.Item("train")
$wsTrain.Activate()
$wsTrain.Range("E3").Value = "0,4"
$excel.ActiveWindow.Zoom = 138

# --- Sheet "pipe-matcher": change the matcher pattern from amazon/BRAND to google/PRODUCT ---
$wsMatcher = $wb.Worksheets.Item("pipe-matcher")
$wsMatcher.Activate()
$wsMatcher.Range("A2").Value = "google"
$wsMatcher.Range("C2").Value = "PRODUCT"
$wsMatcher.Range("C2").Select()

# --- Sheet "config": point the test data file at the tmp dir, bump iteration count ---
$wsConfig = $wb.Worksheets.Item("config")
$wsConfig.Activate()
$wsConfig.Range("B4").Value = "[tmp]/nlp/test_data_04"
$wsConfig.Range("B5").Value = 28
